# Updates Sheet1 (loading_percent.xlsx equivalent) values for "case with 380 kV done"
# Columns B..O, rows 2..25 (A holds the 0..23 index, unchanged)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(9.191640900518481,4.493972008280427,0,23.38996292535205,40.64414866160764,3.634272639075779,0,20.19171980009867,8.007000425280117,9.124514955644681,0,0,17.99883854492609,21.58855934395047),
    @(8.893509616035642,4.263790301143711,0,23.06636014080218,40.46166954390527,3.635843991346709,0,20.28755480560199,8.028862278589747,8.9262423328263,0,0,18.05221474305369,21.67709816221078),
    @(8.706502411890128,4.115036454442852,0,22.87114928959364,40.36050542258614,3.636859314286014,0,20.35010670159825,8.043084394523515,8.803491142628351,0,0,18.08664152811785,21.73596123338912),
    @(8.62941950243922,4.052586378972664,0,22.79257077187284,40.32204656611562,3.637285806598678,0,20.37652987426864,8.049081331312946,8.753283695907772,0,0,18.10108763808136,21.76107802671067),
    @(8.616570679251273,4.042107108440769,0,22.7795843074129,40.31582840102998,3.637357395954019,0,20.38097376423365,8.050089289427513,8.744937652306184,0,0,18.1035116170285,21.76531682111009),
    @(8.705466228362333,4.11420159661355,0,22.87008549353567,40.35997551466182,3.636865014473577,0,20.35045927678188,8.043164455641108,8.802814685149736,0,0,18.08683466386452,21.73629539603318),
    @(9.089733668831185,4.416164271543193,0,23.27771983014616,40.57899167743319,3.63480398273782,0,20.22399403347995,8.014372850892844,9.056400954908451,0,0,18.01690012503799,21.61815242728408),
    @(9.807212276189089,4.948234104331874,0,24.10044947590739,41.09325760072023,3.631161232002963,0,20.00541676941592,7.964231521688066,9.54276705133414,0,0,17.89282820145641,21.42226463696993),
    @(10.3067430795156,5.3012962633359,0,24.71339677445342,41.52044279151271,3.628725547152997,0,19.86275391352352,7.931217611957712,9.889784270186167,0,0,17.80956823868593,21.30028178261236),
    @(10.52700131301659,5.453521842105197,0,24.99287597020901,41.72492605853165,3.627669201394856,0,19.80174426091695,7.917023546496795,10.04470460812652,0,0,17.77339041174816,21.24957714669923),
    @(10.60933225672372,5.509949760732921,0,25.09870158580463,41.80376573306507,3.627276578275172,0,19.77920079092949,7.9117667080224,10.10289226828995,0,0,17.75993380222985,21.23106700821796),
    @(10.59164973650419,5.497851254054962,0,25.07591202682934,41.78672454154695,3.627360808501239,0,19.78403103669574,7.912893613374153,10.09038257179413,0,0,17.76282112484388,21.23502274479038),
    @(10.53379667992386,5.458188636836732,0,25.00158301331441,41.7313843805947,3.62763675211368,0,19.79987837905448,7.916588697463141,10.04950160711747,0,0,17.7722784607561,21.2480404520039),
    @(10.49821787293294,5.433735421069963,0,24.95605049719573,41.69766842688604,3.627806737134727,0,19.80965821851076,7.918867418612684,10.02439712070629,0,0,17.7781029890329,21.25610418179539),
    @(10.29220184957793,5.291178252284502,0,24.69513655087366,41.50727912244486,3.628795618229145,0,19.86681931644109,7.932161779945719,9.879595955441484,0,0,17.81196662161056,21.3036919601653),
    @(10.16397568606162,5.201567966036068,0,24.53516791343556,41.39304485316745,3.629415469940373,0,19.90288198739287,7.940528259018278,9.789972179140054,0,0,17.83317496036539,21.33411302457473),
    @(10.08957008995196,5.149237805319107,0,24.44322487488223,41.32829954171643,3.629776856848722,0,19.92399024776353,7.945418028575121,9.738149294592278,0,0,17.84553327331408,21.35206075203977),
    @(10.06426776635128,5.13138463276045,0,24.41210908404162,41.30654423358429,3.629900052862106,0,19.93119998533213,7.947086958190522,9.720557681684069,0,0,17.84974506738706,21.358214816996),
    @(10.17769374345185,5.211188848950283,0,24.55219067828963,41.40510639192561,3.629348982464629,0,19.89900517556976,7.939629605194083,9.799541513488075,0,0,17.8309007608944,21.33082802493271),
    @(10.55081926422343,5.469871601421106,0,25.02341622335177,41.74760140453246,3.627555500455789,0,19.79520844053754,7.915500158104466,10.06152270766491,0,0,17.76949402052259,21.24419807394434),
    @(10.78837836879904,5.631841141336874,0,25.33130354281135,41.97961173923383,3.626426426338276,0,19.73063305846579,7.900418651380862,10.22993781035975,0,0,17.73077798748693,21.19160688373045),
    @(10.66218682414924,5.546047058222039,0,25.1670185021738,41.85505442748968,3.627025105510059,0,19.76479955000318,7.908405057395103,10.14032518640628,0,0,17.75131213591548,21.21930657897666),
    @(10.17149394293355,5.206841776203446,0,24.54449460054567,41.39965047079481,3.629379025794426,0,19.90075671300713,7.940035638085169,9.795216138250801,0,0,17.83192841131194,21.33231174599161),
    @(9.617593277695551,4.810871016263406,0,23.87597357361232,40.94528474575171,3.632104249659541,0,20.06139923756215,7.977122436108945,9.54276705133414,0,0,17.92500124663307,21.47141380373255)
)

$nRows = $rows.Count
$nCols = $rows[0].Count
$arr = New-Object 'object[,]' $nRows,$nCols
for ($i = 0; $i -lt $nRows; $i++) {
    for ($j = 0; $j -lt $nCols; $j++) {
        $arr[$i,$j] = $rows[$i][$j]
    }
}

$ws.Range("B2:O25").Value = $arr

Write-Output "Updated loading_percent values (case with 380 kV done)"
